$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")

$ws.Range("B2").Value = "2024-06-08"
$ws.Range("C2").Value = "九江·首届萤火之星国风动漫嘉年华"
$ws.Range("D2").Value = "十里大道202号（十里大道与地质路交汇处） 山水国际大酒店"
$ws.Range("E2").Value = "2024.06.08 10:00-06.08 16:00"
$ws.Range("F2").Value = 196
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=85234"
$ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202405/AJrD37gG1715091469262.jpeg"

$ws.Range("B3").Value = "2024-06-08"
$ws.Range("C3").Value = "南昌·CM02动漫游戏博览会"
$ws.Range("D3").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws.Range("E3").Value = "2024.06.08 10:00-06.09 17:00"
$ws.Range("F3").Value = 3176
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85037"
$ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202404/qSrEq0US1713947622923.png"

$ws.Range("B4").Value = "2024-06-08"
$ws.Range("C4").Value = "赣州·次元梦乡2024国风动漫节"
$ws.Range("D4").Value = "廉泉路赣友味餐厅旁 铸谊篮球·羽毛球馆"
$ws.Range("E4").Value = "2024.06.08 10:00-06.08 16:00"
$ws.Range("F4").Value = 232
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85909"
$ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202405/2qcdLboD1715937690231.jpeg"

$ws.Range("B5").Value = "2024-06-09"
$ws.Range("C5").Value = "九江·第四届ACD动漫游戏嘉年华"
$ws.Range("D5").Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"
$ws.Range("E5").Value = "2024.06.09 10:00-06.10 17:00"
$ws.Range("F5").Value = 126
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=85848"
$ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202405/Wsliwm4F1715931131791.jpeg"

$ws.Range("B6").Value = "2024-06-09"
$ws.Range("C6").Value = "信丰·端午节UPUP动漫展"
$ws.Range("D6").Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
$ws.Range("E6").Value = "2024.06.09 10:00-06.09 17:00"
$ws.Range("F6").Value = 199
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84078"
$ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg"

$ws.Range("B7").Value = "2024-06-09"
$ws.Range("C7").Value = "南昌·次元之门动漫游戏嘉年华SP：代号序章"
$ws.Range("D7").Value = "人杰路名实花园北侧约240米 滕王阁游客中心"
$ws.Range("E7").Value = "2024.06.09 10:00-06.10 17:00"
$ws.Range("F7").Value = 1682
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85337"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202405/bpl1MHnz1715237288951.jpeg"

$ws.Range("B8").Value = "2024-06-09"
$ws.Range("C8").Value = "南昌·第三届龙年动漫展——庆端午贺高考专场"
$ws.Range("D8").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws.Range("E8").Value = "2024.06.09 10:00-06.10 18:00"
$ws.Range("F8").Value = 1629
$ws.Range("G8").Value = 55
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=85297"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202405/zBSAcG1V1714936299746.jpeg"

$ws.Range("B9").Value = "2024-06-09"
$ws.Range("C9").Value = "宜春·静卿缤纷仲夏国风动漫文化展"
$ws.Range("D9").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws.Range("E9").Value = "2024.06.09 09:00-06.09 17:00"
$ws.Range("F9").Value = 465
$ws.Range("G9").Value = 45
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85708"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202405/I2WdH04K1715560130445.jpeg"

$ws.Range("B10").Value = "2024-06-10"
$ws.Range("C10").Value = "上饶·ETI动漫节"
$ws.Range("D10").Value = "滨江东路与体育馆路交叉口西100米 力加体育综合运动中心"
$ws.Range("E10").Value = "2024.06.10 10:00-06.10 16:00"
$ws.Range("F10").Value = 369
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=83422"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg"

$ws.Range("B11").Value = "2024-06-10"
$ws.Range("C11").Value = "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）（取消）"
$ws.Range("D11").Value = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
$ws.Range("E11").Value = "2024.06.10 10:00-06.10 17:00"
$ws.Range("F11").Value = 234
$ws.Range("G11").Value = "不可售"
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84575"
$ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg"

$ws.Range("B12").Value = "2024-06-15"
$ws.Range("C12").Value = "上饶·宅舞联萌·随舞动漫派对（免费活动)"
$ws.Range("D12").Value = "春江北大道和吉阳路交汇处 槠溪时光PARK"
$ws.Range("E12").Value = "2024.06.15 08:00-06.15 21:00"
$ws.Range("F12").Value = 28
$ws.Range("G12").Value = 22.33
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85607"
$ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202405/jcZGKqhx1715589649770.jpeg"

$ws.Range("B13").Value = "2024-06-22"
$ws.Range("C13").Value = "景德镇·BM次元盛典运动番only"
$ws.Range("D13").Value = "广场南路金幕影城旁 罗曼园宴会酒店"
$ws.Range("E13").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F13").Value = 190
$ws.Range("G13").Value = 55
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85197"
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png"

$ws.Range("B14").Value = "2024-06-22"
$ws.Range("C14").Value = "萍乡·AU9夏至国漫展"
$ws.Range("D14").Value = "金陵东路18号 萍乡市体育馆"
$ws.Range("E14").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F14").Value = 32
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"

$ws.Range("B15").Value = "2024-06-23"
$ws.Range("C15").Value = "上饶·BM次元盛典运动番only"
$ws.Range("D15").Value = "春江北大道时光PARK内 博悦宴会艺术中心"
$ws.Range("E15").Value = "2024.06.23 10:00-06.23 17:00"
$ws.Range("F15").Value = 229
$ws.Range("G15").Value = 55
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85201"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png"

$ws.Range("B16").Value = "2024-06-29"
$ws.Range("C16").Value = "萍乡·BM次元盛典运动番only"
$ws.Range("D16").Value = "康庄路3号 萍乡梅园国际大酒店"
$ws.Range("E16").Value = "2024.06.29 10:00-06.29 17:00"
$ws.Range("F16").Value = 238
$ws.Range("G16").Value = 55
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=85192"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png"

$ws.Range("B17").Value = "2024-06-30"
$ws.Range("C17").Value = "宜春·BM次元盛典运动番only"
$ws.Range("D17").Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
$ws.Range("E17").Value = "2024.06.30 10:00-06.30 17:00"
$ws.Range("F17").Value = 230
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=84636"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png"

$ws.Range("B18").Value = "2024-07-06"
$ws.Range("C18").Value = "南昌·次元星球动漫游戏展"
$ws.Range("D18").Value = "龙蟠街666号融创茂1层 融创茂"
$ws.Range("E18").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86405"
$ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg"

$ws.Range("B19").Value = "2024-07-06"
$ws.Range("C19").Value = "鹰潭·BM次元盛典运动番only"
$ws.Range("D19").Value = "体育馆东路2号九小隔壁 忆江南•宴会楼"
$ws.Range("E19").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=85997"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png"

$ws.Range("B20").Value = "2024-07-07"
$ws.Range("C20").Value = "赣州·BM次元盛典运动番only"
$ws.Range("D20").Value = "米瑞金路2口0号上客天下1楼 上客天下.老虔州"
$ws.Range("E20").Value = "2024.07.07 10:00-07.07 17:00"
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 55
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=86602"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png"

$ws.Range("B21").Value = "2024-07-12"
$ws.Range("C21").Value = "新余·2024第三届MG动漫嘉年华"
$ws.Range("D21").Value = "仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅"
$ws.Range("E21").Value = "2024.07.12 10:00-07.13 17:30"
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86536"
$ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg"

$ws.Range("B22").Value = "2024-07-13"
$ws.Range("C22").Value = "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华"
$ws.Range("D22").Value = "宜春国际商贸城会展中心 宜春国际商贸城会展中心"
$ws.Range("E22").Value = "2024.07.13 10:00-07.14 17:00"
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 55
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86667"
$ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg"

$ws.Range("B23").Value = "2024-07-14"
$ws.Range("C23").Value = "吉安·COMIC LIFE次元假日05"
$ws.Range("D23").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws.Range("E23").Value = "2024.07.14 09:00-07.14 18:00"
$ws.Range("F23").Value = 376
$ws.Range("G23").Value = 52.1
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=85924"
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg"

$ws.Range("B24").Value = "2024-07-19"
$ws.Range("C24").Value = "赣州·第四届赣州半夏动漫展"
$ws.Range("D24").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
$ws.Range("E24").Value = "2024.07.19 10:00-07.21 17:00"
$ws.Range("F24").Value = 205
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=86587"
$ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg"

$ws.Range("B25").Value = "2024-07-20"
$ws.Range("C25").Value = "南昌·漫拥动漫嘉年华Pro-追光启航"
$ws.Range("D25").Value = "小蓝南路420号 洪州体育馆"
$ws.Range("E25").Value = "2024.07.20 09:00-07.21 17:00"
$ws.Range("F25").Value = 102
$ws.Range("G25").Value = 52.5
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=85796"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png"

$ws.Range("B26").Value = "2024-07-21"
$ws.Range("C26").Value = "乐平·CY境界次元动漫夏时庆"
$ws.Range("D26").Value = "翥山西路182号 佳佳基大酒店"
$ws.Range("E26").Value = "2024.07.21 10:00-07.21 17:00"
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=86768"
$ws.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png"

$ws.Range("B27").Value = "2024-07-21"
$ws.Range("C27").Value = "九江·SXD动漫嘉年华"
$ws.Range("D27").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws.Range("E27").Value = "2024.07.21 10:00-07.21 17:30"
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

$ws.Range("B28").Value = "2024-07-21"
$ws.Range("C28").Value = "萍乡·NL14动漫游戏展·夏日狂想曲"
$ws.Range("D28").Value = "公园南路168号(近工行城北分理处) 梅生嘉华酒店"
$ws.Range("E28").Value = "2024.07.21 10:00-07.21 17:00"
$ws.Range("F28").Value = 23
$ws.Range("G28").Value = 40
$ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86658"
$ws.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg"

$ws.Range("B29").Value = "2024-07-26"
$ws.Range("C29").Value = "南昌·萌卡动漫展"
$ws.Range("D29").Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
$ws.Range("E29").Value = "2024.07.26 09:00-07.28 17:00"
$ws.Range("F29").Value = 219
$ws.Range("G29").Value = 65
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=86776"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg"

$ws.Range("B30").Value = "2024-07-27"
$ws.Range("C30").Value = "江西·次元星河动漫游戏嘉年华"
$ws.Range("D30").Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws.Range("E30").Value = "2024.07.27 10:00-07.28 17:00"
$ws.Range("F30").Value = 2152
$ws.Range("G30").Value = 69
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png"

$ws.Range("B31").Value = "2024-07-27"
$ws.Range("C31").Value = "赣州·马娘only"
$ws.Range("D31").Value = "火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)"
$ws.Range("E31").Value = "2024.07.27 09:00-07.27 17:00"
$ws.Range("F31").Value = 7
$ws.Range("G31").Value = 60
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=86772"
$ws.Range("I31").Value = "//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png"

$ws.Range("B32").Value = "2024-07-28"
$ws.Range("C32").Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws.Range("D32").Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws.Range("E32").Value = "2024.07.28 11:00-07.28 17:00"
$ws.Range("F32").Value = 52
$ws.Range("G32").Value = 56
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"

$ws.Range("B33").Value = "2024-08-03"
$ws.Range("C33").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws.Range("D33").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws.Range("E33").Value = "2024.08.03 09:00-08.04 17:30"
$ws.Range("F33").Value = 467
$ws.Range("G33").Value = 64
$ws.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"

$ws.Range("B34").Value = "2024-08-03"
$ws.Range("C34").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws.Range("D34").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws.Range("E34").Value = "2024.08.03 09:00-08.04 17:00"
$ws.Range("F34").Value = 321
$ws.Range("G34").Value = 55
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"

$ws.Range("B35").Value = "2024-08-03"
$ws.Range("C35").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws.Range("D35").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws.Range("E35").Value = "2024.08.03 08:30-08.03 17:00"
$ws.Range("F35").Value = 570
$ws.Range("G35").Value = "已售罄"
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"

$ws.Range("B36").Value = "2024-08-03"
$ws.Range("C36").Value = "樟树·第二届静卿国风动漫文化展览会"
$ws.Range("D36").Value = "杏佛路89号 樟树银河国际酒店"
$ws.Range("E36").Value = "2024.08.03 09:00-08.03 17:00"
$ws.Range("F36").Value = 425
$ws.Range("G36").Value = 45
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"

$ws.Range("B37").Value = "2024-08-04"
$ws.Range("C37").Value = "九江·第一届异次元动漫嘉年华"
$ws.Range("D37").Value = "长虹西大道兴城广场99号 九江半岛宾馆"
$ws.Range("E37").Value = "2024.08.04 08:00-08.04 17:00"
$ws.Range("F37").Value = 226
$ws.Range("G37").Value = 45
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=84407"
$ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg"

$ws.Range("B38").Value = "2024-08-06"
$ws.Range("C38").Value = "南昌·第一届异次元动漫嘉年华"
$ws.Range("D38").Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws.Range("E38").Value = "2024.08.06 08:00-08.06 17:00"
$ws.Range("F38").Value = 343
$ws.Range("G38").Value = 55
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"

$ws.Range("B39").Value = "2024-08-06"
$ws.Range("C39").Value = "宜春·第三十五届静卿国风动漫文化展览会"
$ws.Range("D39").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws.Range("E39").Value = "2024.08.06 09:00-08.06 17:00"
$ws.Range("F39").Value = 412
$ws.Range("G39").Value = 45
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=86684"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg"

$ws.Range("B40").Value = "2024-08-08"
$ws.Range("C40").Value = "赣州·第二届异次元动漫嘉年华"
$ws.Range("D40").Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
$ws.Range("E40").Value = "2024.08.08 08:00-08.08 17:00"
$ws.Range("F40").Value = 513
$ws.Range("G40").Value = 45
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=84184"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg"

$ws.Range("B41").Value = "2024-08-10"
$ws.Range("C41").Value = "高安·第二届静卿国风动漫文化展览会"
$ws.Range("D41").Value = "华林中路606号 高安华鼎国际大酒店"
$ws.Range("E41").Value = "2024.08.10 09:00-08.10 17:00"
$ws.Range("F41").Value = 414
$ws.Range("G41").Value = 45
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=86682"
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg"

$ws.Range("A42:I45").EntireRow.Delete()

$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("B2").Value = "2024-06-08"
$ws.Range("C2").Value = "九江·首届萤火之星国风动漫嘉年华"
$ws.Range("D2").Value = "十里大道202号（十里大道与地质路交汇处） 山水国际大酒店"
$ws.Range("E2").Value = "2024.06.08 10:00-06.08 16:00"
$ws.Range("F2").Value = 196
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=85234"
$ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202405/AJrD37gG1715091469262.jpeg"

$ws.Range("B3").Value = "2024-06-08"
$ws.Range("C3").Value = "南昌·CM02动漫游戏博览会"
$ws.Range("D3").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
$ws.Range("E3").Value = "2024.06.08 10:00-06.09 17:00"
$ws.Range("F3").Value = 3176
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=85037"
$ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202404/qSrEq0US1713947622923.png"

$ws.Range("B4").Value = "2024-06-08"
$ws.Range("C4").Value = "赣州·次元梦乡2024国风动漫节"
$ws.Range("D4").Value = "廉泉路赣友味餐厅旁 铸谊篮球·羽毛球馆"
$ws.Range("E4").Value = "2024.06.08 10:00-06.08 16:00"
$ws.Range("F4").Value = 232
$ws.Range("G4").Value = 55
$ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=85909"
$ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202405/2qcdLboD1715937690231.jpeg"

$ws.Range("B5").Value = "2024-06-09"
$ws.Range("C5").Value = "九江·第四届ACD动漫游戏嘉年华"
$ws.Range("D5").Value = "九瑞大道与重庆路交汇处西南角 九江国际会展中心"
$ws.Range("E5").Value = "2024.06.09 10:00-06.10 17:00"
$ws.Range("F5").Value = 126
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=85848"
$ws.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202405/Wsliwm4F1715931131791.jpeg"

$ws.Range("B6").Value = "2024-06-09"
$ws.Range("C6").Value = "信丰·端午节UPUP动漫展"
$ws.Range("D6").Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
$ws.Range("E6").Value = "2024.06.09 10:00-06.09 17:00"
$ws.Range("F6").Value = 199
$ws.Range("G6").Value = 48
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84078"
$ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg"

$ws.Range("B7").Value = "2024-06-09"
$ws.Range("C7").Value = "南昌·次元之门动漫游戏嘉年华SP：代号序章"
$ws.Range("D7").Value = "人杰路名实花园北侧约240米 滕王阁游客中心"
$ws.Range("E7").Value = "2024.06.09 10:00-06.10 17:00"
$ws.Range("F7").Value = 1682
$ws.Range("G7").Value = 60
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=85337"
$ws.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202405/bpl1MHnz1715237288951.jpeg"

$ws.Range("B8").Value = "2024-06-09"
$ws.Range("C8").Value = "南昌·第三届龙年动漫展——庆端午贺高考专场"
$ws.Range("D8").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws.Range("E8").Value = "2024.06.09 10:00-06.10 18:00"
$ws.Range("F8").Value = 1629
$ws.Range("G8").Value = 55
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=85297"
$ws.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202405/zBSAcG1V1714936299746.jpeg"

$ws.Range("B9").Value = "2024-06-09"
$ws.Range("C9").Value = "宜春·静卿缤纷仲夏国风动漫文化展"
$ws.Range("D9").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws.Range("E9").Value = "2024.06.09 09:00-06.09 17:00"
$ws.Range("F9").Value = 465
$ws.Range("G9").Value = 45
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=85708"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202405/I2WdH04K1715560130445.jpeg"

$ws.Range("B10").Value = "2024-06-10"
$ws.Range("C10").Value = "上饶·ETI动漫节"
$ws.Range("D10").Value = "滨江东路与体育馆路交叉口西100米 力加体育综合运动中心"
$ws.Range("E10").Value = "2024.06.10 10:00-06.10 16:00"
$ws.Range("F10").Value = 369
$ws.Range("G10").Value = 55
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=83422"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202404/N6VdMOuL1713257425864.jpeg"

$ws.Range("B11").Value = "2024-06-10"
$ws.Range("C11").Value = "南昌·LY-COSPLAY大会X运动番PRO2.0（非ONLY）（取消）"
$ws.Range("D11").Value = "青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK"
$ws.Range("E11").Value = "2024.06.10 10:00-06.10 17:00"
$ws.Range("F11").Value = 234
$ws.Range("G11").Value = "不可售"
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84575"
$ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202404/ScwkijwU1713428452963.jpeg"

$ws.Range("B12").Value = "2024-06-15"
$ws.Range("C12").Value = "上饶·宅舞联萌·随舞动漫派对（免费活动)"
$ws.Range("D12").Value = "春江北大道和吉阳路交汇处 槠溪时光PARK"
$ws.Range("E12").Value = "2024.06.15 08:00-06.15 21:00"
$ws.Range("F12").Value = 28
$ws.Range("G12").Value = 22.33
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85607"
$ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202405/jcZGKqhx1715589649770.jpeg"

$ws.Range("B13").Value = "2024-06-22"
$ws.Range("C13").Value = "景德镇·BM次元盛典运动番only"
$ws.Range("D13").Value = "广场南路金幕影城旁 罗曼园宴会酒店"
$ws.Range("E13").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F13").Value = 190
$ws.Range("G13").Value = 55
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=85197"
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202404/Z6eXz0su1714292081978.png"

$ws.Range("B14").Value = "2024-06-22"
$ws.Range("C14").Value = "萍乡·AU9夏至国漫展"
$ws.Range("D14").Value = "金陵东路18号 萍乡市体育馆"
$ws.Range("E14").Value = "2024.06.22 10:00-06.22 17:00"
$ws.Range("F14").Value = 32
$ws.Range("G14").Value = 45
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=86453"
$ws.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202405/iFDRERFO1716547195192.jpeg"

$ws.Range("B15").Value = "2024-06-23"
$ws.Range("C15").Value = "上饶·BM次元盛典运动番only"
$ws.Range("D15").Value = "春江北大道时光PARK内 博悦宴会艺术中心"
$ws.Range("E15").Value = "2024.06.23 10:00-06.23 17:00"
$ws.Range("F15").Value = 229
$ws.Range("G15").Value = 55
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=85201"
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202404/30dgkbjT1714293499693.png"

$ws.Range("B16").Value = "2024-06-29"
$ws.Range("C16").Value = "萍乡·BM次元盛典运动番only"
$ws.Range("D16").Value = "康庄路3号 萍乡梅园国际大酒店"
$ws.Range("E16").Value = "2024.06.29 10:00-06.29 17:00"
$ws.Range("F16").Value = 238
$ws.Range("G16").Value = 55
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=85192"
$ws.Range("I16").Value = "//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png"

$ws.Range("B17").Value = "2024-06-30"
$ws.Range("C17").Value = "宜春·BM次元盛典运动番only"
$ws.Range("D17").Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
$ws.Range("E17").Value = "2024.06.30 10:00-06.30 17:00"
$ws.Range("F17").Value = 230
$ws.Range("G17").Value = 55
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=84636"
$ws.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png"

$ws.Range("B18").Value = "2024-07-06"
$ws.Range("C18").Value = "南昌·次元星球动漫游戏展"
$ws.Range("D18").Value = "龙蟠街666号融创茂1层 融创茂"
$ws.Range("E18").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 50
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=86405"
$ws.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg"

$ws.Range("B19").Value = "2024-07-06"
$ws.Range("C19").Value = "鹰潭·BM次元盛典运动番only"
$ws.Range("D19").Value = "体育馆东路2号九小隔壁 忆江南•宴会楼"
$ws.Range("E19").Value = "2024.07.06 10:00-07.06 17:00"
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=85997"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png"

$ws.Range("B20").Value = "2024-07-07"
$ws.Range("C20").Value = "赣州·BM次元盛典运动番only"
$ws.Range("D20").Value = "米瑞金路2口0号上客天下1楼 上客天下.老虔州"
$ws.Range("E20").Value = "2024.07.07 10:00-07.07 17:00"
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 55
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=86602"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png"

$ws.Range("B21").Value = "2024-07-12"
$ws.Range("C21").Value = "新余·2024第三届MG动漫嘉年华"
$ws.Range("D21").Value = "仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅"
$ws.Range("E21").Value = "2024.07.12 10:00-07.13 17:30"
$ws.Range("F21").Value = 51
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=86536"
$ws.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg"

$ws.Range("B22").Value = "2024-07-13"
$ws.Range("C22").Value = "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华"
$ws.Range("D22").Value = "宜春国际商贸城会展中心 宜春国际商贸城会展中心"
$ws.Range("E22").Value = "2024.07.13 10:00-07.14 17:00"
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 55
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86667"
$ws.Range("I22").Value = "//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg"

$ws.Range("B23").Value = "2024-07-14"
$ws.Range("C23").Value = "吉安·COMIC LIFE次元假日05"
$ws.Range("D23").Value = "东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心"
$ws.Range("E23").Value = "2024.07.14 09:00-07.14 18:00"
$ws.Range("F23").Value = 376
$ws.Range("G23").Value = 52.1
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=85924"
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg"

$ws.Range("B24").Value = "2024-07-19"
$ws.Range("C24").Value = "赣州·第四届赣州半夏动漫展"
$ws.Range("D24").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
$ws.Range("E24").Value = "2024.07.19 10:00-07.21 17:00"
$ws.Range("F24").Value = 205
$ws.Range("G24").Value = 55
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=86587"
$ws.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg"

$ws.Range("B25").Value = "2024-07-20"
$ws.Range("C25").Value = "南昌·漫拥动漫嘉年华Pro-追光启航"
$ws.Range("D25").Value = "小蓝南路420号 洪州体育馆"
$ws.Range("E25").Value = "2024.07.20 09:00-07.21 17:00"
$ws.Range("F25").Value = 102
$ws.Range("G25").Value = 52.5
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=85796"
$ws.Range("I25").Value = "//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png"

$ws.Range("B26").Value = "2024-07-21"
$ws.Range("C26").Value = "乐平·CY境界次元动漫夏时庆"
$ws.Range("D26").Value = "翥山西路182号 佳佳基大酒店"
$ws.Range("E26").Value = "2024.07.21 10:00-07.21 17:00"
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=86768"
$ws.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png"

$ws.Range("B27").Value = "2024-07-21"
$ws.Range("C27").Value = "九江·SXD动漫嘉年华"
$ws.Range("D27").Value = "湓浦街道大中路339号 百嘉洲际酒店"
$ws.Range("E27").Value = "2024.07.21 10:00-07.21 17:30"
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=86832"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg"

$ws.Range("B28").Value = "2024-07-21"
$ws.Range("C28").Value = "萍乡·NL14动漫游戏展·夏日狂想曲"
$ws.Range("D28").Value = "公园南路168号(近工行城北分理处) 梅生嘉华酒店"
$ws.Range("E28").Value = "2024.07.21 10:00-07.21 17:00"
$ws.Range("F28").Value = 23
$ws.Range("G28").Value = 40
$ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=86658"
$ws.Range("I28").Value = "//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg"

$ws.Range("B29").Value = "2024-07-26"
$ws.Range("C29").Value = "南昌·萌卡动漫展"
$ws.Range("D29").Value = "八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆"
$ws.Range("E29").Value = "2024.07.26 09:00-07.28 17:00"
$ws.Range("F29").Value = 219
$ws.Range("G29").Value = 65
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=86776"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg"

$ws.Range("B30").Value = "2024-07-27"
$ws.Range("C30").Value = "江西·次元星河动漫游戏嘉年华"
$ws.Range("D30").Value = "九龙大道1177号 南昌绿地国际博览中心"
$ws.Range("E30").Value = "2024.07.27 10:00-07.28 17:00"
$ws.Range("F30").Value = 2152
$ws.Range("G30").Value = 69
$ws.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=85493"
$ws.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png"

$ws.Range("B31").Value = "2024-07-27"
$ws.Range("C31").Value = "赣州·马娘only"
$ws.Range("D31").Value = "火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)"
$ws.Range("E31").Value = "2024.07.27 09:00-07.27 17:00"
$ws.Range("F31").Value = 7
$ws.Range("G31").Value = 60
$ws.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=86772"
$ws.Range("I31").Value = "//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png"

$ws.Range("B32").Value = "2024-07-28"
$ws.Range("C32").Value = "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会"
$ws.Range("D32").Value = "兴国路恒大帝景西门 江西长庚控股有限公司"
$ws.Range("E32").Value = "2024.07.28 11:00-07.28 17:00"
$ws.Range("F32").Value = 52
$ws.Range("G32").Value = 56
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=85688"
$ws.Range("I32").Value = "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png"

$ws.Range("B33").Value = "2024-08-03"
$ws.Range("C33").Value = "南昌·幻梦境国际动漫游戏嘉年华1th"
$ws.Range("D33").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws.Range("E33").Value = "2024.08.03 09:00-08.04 17:30"
$ws.Range("F33").Value = 467
$ws.Range("G33").Value = 64
$ws.Range("H33").Value = "https://show.bilibili.com/platform/detail.html?id=83980"
$ws.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg"

$ws.Range("B34").Value = "2024-08-03"
$ws.Range("C34").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会"
$ws.Range("D34").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws.Range("E34").Value = "2024.08.03 09:00-08.04 17:00"
$ws.Range("F34").Value = 321
$ws.Range("G34").Value = 55
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=86341"
$ws.Range("I34").Value = "//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png"

$ws.Range("B35").Value = "2024-08-03"
$ws.Range("C35").Value = "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票"
$ws.Range("D35").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
$ws.Range("E35").Value = "2024.08.03 08:30-08.03 17:00"
$ws.Range("F35").Value = 570
$ws.Range("G35").Value = "已售罄"
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=85981"
$ws.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png"

$ws.Range("B36").Value = "2024-08-03"
$ws.Range("C36").Value = "樟树·第二届静卿国风动漫文化展览会"
$ws.Range("D36").Value = "杏佛路89号 樟树银河国际酒店"
$ws.Range("E36").Value = "2024.08.03 09:00-08.03 17:00"
$ws.Range("F36").Value = 425
$ws.Range("G36").Value = 45
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=86683"
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg"

$ws.Range("B37").Value = "2024-08-04"
$ws.Range("C37").Value = "九江·第一届异次元动漫嘉年华"
$ws.Range("D37").Value = "长虹西大道兴城广场99号 九江半岛宾馆"
$ws.Range("E37").Value = "2024.08.04 08:00-08.04 17:00"
$ws.Range("F37").Value = 226
$ws.Range("G37").Value = 45
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=84407"
$ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202404/e7k26XLV1713262153782.jpeg"

$ws.Range("B38").Value = "2024-08-06"
$ws.Range("C38").Value = "南昌·第一届异次元动漫嘉年华"
$ws.Range("D38").Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws.Range("E38").Value = "2024.08.06 08:00-08.06 17:00"
$ws.Range("F38").Value = 343
$ws.Range("G38").Value = 55
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=84102"
$ws.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg"

$ws.Range("B39").Value = "2024-08-06"
$ws.Range("C39").Value = "宜春·第三十五届静卿国风动漫文化展览会"
$ws.Range("D39").Value = "宜阳大道19号(交通银行旁) 宜春安缦文华酒店"
$ws.Range("E39").Value = "2024.08.06 09:00-08.06 17:00"
$ws.Range("F39").Value = 412
$ws.Range("G39").Value = 45
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=86684"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg"

$ws.Range("B40").Value = "2024-08-08"
$ws.Range("C40").Value = "赣州·第二届异次元动漫嘉年华"
$ws.Range("D40").Value = "金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆"
$ws.Range("E40").Value = "2024.08.08 08:00-08.08 17:00"
$ws.Range("F40").Value = 513
$ws.Range("G40").Value = 45
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=84184"
$ws.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg"

$ws.Range("B41").Value = "2024-08-10"
$ws.Range("C41").Value = "高安·第二届静卿国风动漫文化展览会"
$ws.Range("D41").Value = "华林中路606号 高安华鼎国际大酒店"
$ws.Range("E41").Value = "2024.08.10 09:00-08.10 17:00"
$ws.Range("F41").Value = 414
$ws.Range("G41").Value = 45
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=86682"
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg"

$ws.Range("A42:I45").EntireRow.Delete()
